$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Add a new worksheet right after Sheet1 and name it "Sheet2"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Find the last used row on Sheet1 (column A)
$lastRow = $ws1.Cells.Item($ws1.Rows.Count, 1).End(-4162).Row

# Build a multi-area range: header row + every data row whose Type
# (column C) is "WI3" (Process Vendor Invoice) or "WI5" (Calculate
# Client Security Hash)
$srcRange = $ws1.Range("A1:E1")
for ($r = 2; $r -le $lastRow; $r++) {
    $type = $ws1.Cells.Item($r, 3).Value2
    if ($type -eq "WI3" -or $type -eq "WI5") {
        $rowRange = $ws1.Range("A" + $r + ":E" + $r)
        $srcRange = $excel.Union($srcRange, $rowRange)
    }
}

# Copy the filtered rows (values only, preserving original cell typing)
# into the new sheet, starting at A1
$srcRange.Copy()
$ws2.Range("A1").PasteSpecial(-4163)
